$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" column (C) for all data rows (2-10): 46063 -> 46064
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 46064
}

# Capture current (pre-reorder) values for rows 4, 5, 7, 8, 9 in columns A, B, G
$row4A = $ws.Cells.Item(4, 1).Value2
$row4B = $ws.Cells.Item(4, 2).Value2
$row4G = $ws.Cells.Item(4, 7).Value2

$row5A = $ws.Cells.Item(5, 1).Value2
$row5B = $ws.Cells.Item(5, 2).Value2
$row5G = $ws.Cells.Item(5, 7).Value2

$row7A = $ws.Cells.Item(7, 1).Value2
$row7B = $ws.Cells.Item(7, 2).Value2
$row7G = $ws.Cells.Item(7, 7).Value2

$row8A = $ws.Cells.Item(8, 1).Value2
$row8B = $ws.Cells.Item(8, 2).Value2
$row8G = $ws.Cells.Item(8, 7).Value2

$row9A = $ws.Cells.Item(9, 1).Value2
$row9B = $ws.Cells.Item(9, 2).Value2
$row9G = $ws.Cells.Item(9, 7).Value2

# Apply the new order:
# new row4 <- old row7
$ws.Cells.Item(4, 1).Value2 = $row7A
$ws.Cells.Item(4, 2).Value2 = $row7B
$ws.Cells.Item(4, 7).Value2 = $row7G

# new row5 <- old row8
$ws.Cells.Item(5, 1).Value2 = $row8A
$ws.Cells.Item(5, 2).Value2 = $row8B
$ws.Cells.Item(5, 7).Value2 = $row8G

# new row7 <- old row4
$ws.Cells.Item(7, 1).Value2 = $row4A
$ws.Cells.Item(7, 2).Value2 = $row4B
$ws.Cells.Item(7, 7).Value2 = $row4G

# new row8 <- old row9
$ws.Cells.Item(8, 1).Value2 = $row9A
$ws.Cells.Item(8, 2).Value2 = $row9B
$ws.Cells.Item(8, 7).Value2 = $row9G

# new row9 <- old row5
$ws.Cells.Item(9, 1).Value2 = $row5A
$ws.Cells.Item(9, 2).Value2 = $row5B
$ws.Cells.Item(9, 7).Value2 = $row5G
